# Update cryptocurrency price/volume data (cryptos list refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.048.46'
$ws.Range('E2').Value = '  +0.12%  '
$ws.Range('D3').Value = '2.301.02'
$ws.Range('E3').Value = '  +0.29%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '300.04'
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '97.82'
$ws.Range('E6').Value = '  -1.40%  '
$ws.Range('E7').Value = '  +3.50%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.517'
$ws.Range('E9').Value = '  +1.53%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '36.36'
$ws.Range('E10').Value = '  +0.54%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0793'
$ws.Range('E11').Value = '  +0.58%  '
$ws.Range('E12').Value = '  +0.48%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '17.77'
$ws.Range('E13').Value = '  -1.70%  '
$ws.Range('E14').Value = '  -0.65%  '
$ws.Range('D15').Value = '2.658.09'
$ws.Range('E15').Value = '  +0.25%  '
$ws.Range('D16').Value = '2.275.12'
$ws.Range('E16').Value = '  -0.87%  '
$ws.Range('E17').Value = '  -1.16%  '
$ws.Range('D18').Value = '42.940.59'
$ws.Range('E18').Value = '  +0.12%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '13.00'
$ws.Range('E19').Value = '  +3.69%  '
$ws.Range('D20').Value = '0.0₃0912'
$ws.Range('E20').Value = '  +1.31%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.14'
$ws.Range('E21').Value = '  +0.68%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '68.29'
$ws.Range('E22').Value = '  +0.85%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '238.07'
$ws.Range('E23').Value = '  +0.95%  '
$ws.Range('E24').Value = '  -0.69%  '
$ws.Range('E25').Value = '  -0.07%  '
$ws.Range('E26').Value = '  -0.48%  '
$ws.Range('E27').Value = '  -0.04%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '24.97'
$ws.Range('E28').Value = '  +0.14%  '
$ws.Range('E29').Value = '  -12.79%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '163.35'
$ws.Range('E31').Value = '  -2.32%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '33.10'
$ws.Range('E32').Value = '  -4.04%  '
$ws.Range('E33').Value = '  +0.06%  '
$ws.Range('E34').Value = '  +2.35%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '18.19'
$ws.Range('E35').Value = '  +3.27%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '4.75'
$ws.Range('E36').Value = '  +3.03%  '
$ws.Range('E37').Value = '  +0.43%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.0698'
$ws.Range('E38').Value = '  +1.52%  '
$ws.Range('E39').Value = '  +0.63%  '
$ws.Range('E40').Value = '  -0.21%  '
$ws.Range('E41').Value = '  +1.89%  '
$ws.Range('E42').Value = '  -1.32%  '
$ws.Range('D43').Value = '2.014.11'
$ws.Range('E43').Value = '  +2.30%  '
$ws.Range('E44').Value = '  -1.62%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '2.23'
$ws.Range('E45').Value = '  -2.91%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '10.31'
$ws.Range('E46').Value = '  +1.39%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '17.50'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.85'
$ws.Range('E48').Value = '  -1.53%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '54.32'
$ws.Range('E49').Value = '  -1.95%  '
$ws.Range('D50').Value = '2.530.07'
$ws.Range('E50').Value = '  +0.42%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.54'
$ws.Range('E51').Value = '  -0.52%  '

$wb.Save()
